$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 572, shifting existing rows 572:627 down to 573:628
$ws.Rows.Item(572).Insert()

# Populate the newly inserted row 572 with its data
$ws.Cells.Item(572, 1).Value = 5
$ws.Cells.Item(572, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(572, 3).Value = "Maule"
$ws.Cells.Item(572, 4).Value = 45132
$ws.Cells.Item(572, 5).Value = 7
$ws.Cells.Item(572, 6).Value = 100114014
$ws.Cells.Item(572, 7).Value = "Betarraga"
$ws.Cells.Item(572, 8).Value = "Sin especificar"
$ws.Cells.Item(572, 9).Value = "Primera"
$ws.Cells.Item(572, 10).Value = 5000
$ws.Cells.Item(572, 11).Value = 600
$ws.Cells.Item(572, 12).Value = 600
$ws.Cells.Item(572, 13).Value = 600
$ws.Cells.Item(572, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(572, 15).Value = "Región del Maule"
$ws.Cells.Item(572, 16).Value = 120
$ws.Cells.Item(572, 17).Value = 5
$ws.Cells.Item(572, 18).Value = "Hortaliza"
